# Auto-generated edit script for Brynhildr_Profits workbook update
# Updates market price / profit calculation columns (H-N) across multiple leve tables
# to reflect refreshed market data from the scheduled runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 157296.14
$ws.Range("I6").Value = 166845.67
$ws.Range("K6").Value = 500537.01
$ws.Range("M6").Value = -500425.01
$ws.Range("H9").Value = 205.05556
$ws.Range("I9").Value = 103.35714
$ws.Range("K9").Value = 103.35714
$ws.Range("M9").Value = 65.64286
$ws.Range("H12").Value = 14395.714
$ws.Range("I12").Value = 16761.666
$ws.Range("K12").Value = 16761.666
$ws.Range("M12").Value = -16591.666
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("H28").Value = 740
$ws.Range("I28").Value = 816.6667
$ws.Range("J28").Value = 280
$ws.Range("K28").Value = 816.6667
$ws.Range("L28").Value = 280
$ws.Range("M28").Value = -331.6667
$ws.Range("N28").Value = -1250
$ws.Range("H29").Value = 1733.3334
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("H33").Value = 108.5
$ws.Range("I33").Value = 103.26667
$ws.Range("K33").Value = 103.26667
$ws.Range("M33").Value = 125.73333
$ws.Range("H38").Value = 1209.1538
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("H98").Value = 1314.625
$ws.Range("I98").Value = 1314.625
$ws.Range("K98").Value = 1314.625
$ws.Range("M98").Value = 183.375
$ws.Range("H101").Value = 299
$ws.Range("I101").Value = 299
$ws.Range("K101").Value = 897
$ws.Range("M101").Value = 725
$ws.Range("H122").Value = 1314.625
$ws.Range("I122").Value = 1314.625
$ws.Range("K122").Value = 3943.875
$ws.Range("M122").Value = -1493.875
$ws.Range("H138").Value = 4500.773
$ws.Range("J138").Value = 7588.8
$ws.Range("L138").Value = 22766.4
$ws.Range("N138").Value = -33046.4
$ws.Range("M21").ClearContents()
$ws.Range("M23").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("N38").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 778.7931
$ws.Range("I2").Value = 710.35
$ws.Range("J2").Value = 930.8889
$ws.Range("K2").Value = 710.35
$ws.Range("L2").Value = 930.8889
$ws.Range("M2").Value = -597.35
$ws.Range("N2").Value = -1156.8889
$ws.Range("H25").Value = 4751.857
$ws.Range("I25").Value = 815.75
$ws.Range("K25").Value = 815.75
$ws.Range("M25").Value = -413.75
$ws.Range("H32").Value = 170009.94
$ws.Range("I32").Value = 291389.3
$ws.Range("K32").Value = 291389.3
$ws.Range("M32").Value = -291102.3
$ws.Range("H45").Value = 1834.4286
$ws.Range("J45").Value = 1997.6
$ws.Range("L45").Value = 1997.6
$ws.Range("N45").Value = -2751.6
$ws.Range("H102").Value = 2551.0833
$ws.Range("I102").Value = 2801.3
$ws.Range("J102").Value = 1300
$ws.Range("K102").Value = 2801.3
$ws.Range("L102").Value = 1300
$ws.Range("M102").Value = -1179.3
$ws.Range("N102").Value = -4544
$ws.Range("H116").Value = 778.7931
$ws.Range("I116").Value = 710.35
$ws.Range("J116").Value = 930.8889
$ws.Range("K116").Value = 710.35
$ws.Range("L116").Value = 930.8889
$ws.Range("M116").Value = 1583.65
$ws.Range("N116").Value = -5518.8889

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 778.7931
$ws.Range("I3").Value = 710.35
$ws.Range("J3").Value = 930.8889
$ws.Range("K3").Value = 710.35
$ws.Range("L3").Value = 930.8889
$ws.Range("M3").Value = -596.35
$ws.Range("N3").Value = -1158.8889
$ws.Range("H20").Value = 4153.1
$ws.Range("I20").Value = 4960.933
$ws.Range("J20").Value = 1729.6
$ws.Range("K20").Value = 4960.933
$ws.Range("L20").Value = 1729.6
$ws.Range("M20").Value = -4713.933
$ws.Range("N20").Value = -2223.6
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("H134").Value = 2783369.2
$ws.Range("I134").Value = 3932.4614
$ws.Range("K134").Value = 11797.3842
$ws.Range("M134").Value = -9262.3842
$ws.Range("N46").ClearContents()
$ws.Range("N76").ClearContents()
$ws.Range("N79").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4276856
$ws.Range("I31").Value = 5053739
$ws.Range("K31").Value = 5053739
$ws.Range("M31").Value = -5053444
$ws.Range("H34").Value = 4276856
$ws.Range("I34").Value = 5053739
$ws.Range("K34").Value = 5053739
$ws.Range("M34").Value = -5053537
$ws.Range("H134").Value = 1386.762
$ws.Range("I134").Value = 1323.0244
$ws.Range("K134").Value = 3969.0732
$ws.Range("M134").Value = -1434.0732

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1457.6666
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("H56").Value = 8082.8335
$ws.Range("I56").Value = 8082.8335
$ws.Range("K56").Value = 8082.8335
$ws.Range("M56").Value = -7552.8335
$ws.Range("H70").Value = 2122
$ws.Range("I70").Value = 696
$ws.Range("J70").Value = 2205.8823
$ws.Range("K70").Value = 2088
$ws.Range("L70").Value = 6617.646900000001
$ws.Range("M70").Value = -1773
$ws.Range("N70").Value = -7247.646900000001
$ws.Range("H73").Value = 2122
$ws.Range("I73").Value = 696
$ws.Range("J73").Value = 2205.8823
$ws.Range("K73").Value = 2088
$ws.Range("L73").Value = 6617.646900000001
$ws.Range("M73").Value = -996
$ws.Range("N73").Value = -8801.6469
$ws.Range("N31").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9533.878000000001
$ws.Range("I132").Value = 8187.6587
$ws.Range("J132").Value = 16433.25
$ws.Range("K132").Value = 24562.9761
$ws.Range("L132").Value = 49299.75
$ws.Range("M132").Value = -22032.9761
$ws.Range("N132").Value = -54359.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 22950
$ws.Range("J20").Value = 22950
$ws.Range("L20").Value = 22950
$ws.Range("N20").Value = -23402
$ws.Range("H22").Value = 2188.9167
$ws.Range("J22").Value = 4249.25
$ws.Range("L22").Value = 4249.25
$ws.Range("N22").Value = -4839.25
$ws.Range("H27").Value = 2188.9167
$ws.Range("J27").Value = 4249.25
$ws.Range("L27").Value = 4249.25
$ws.Range("N27").Value = -4463.25
$ws.Range("H38").Value = 10665.667
$ws.Range("J38").Value = 10665.667
$ws.Range("L38").Value = 10665.667
$ws.Range("N38").Value = -11485.667
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H99").Value = 16564.5
$ws.Range("I99").Value = 16564.5
$ws.Range("K99").Value = 16564.5
$ws.Range("M99").Value = -13569.5
$ws.Range("H122").Value = 3181.7
$ws.Range("I122").Value = 2791.1875
$ws.Range("K122").Value = 8373.5625
$ws.Range("M122").Value = -5923.5625
$ws.Range("N92").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 16828.846
$ws.Range("I51").Value = 12504.889
$ws.Range("J51").Value = 26557.75
$ws.Range("K51").Value = 12504.889
$ws.Range("L51").Value = 26557.75
$ws.Range("M51").Value = -11994.889
$ws.Range("N51").Value = -27577.75
$ws.Range("H52").Value = 18947.166
$ws.Range("I52").Value = 9847.333000000001
$ws.Range("K52").Value = 9847.333000000001
$ws.Range("M52").Value = -9621.333000000001
$ws.Range("H122").Value = 48261.24
$ws.Range("I122").Value = 3487.6086
$ws.Range("J122").Value = 563158
$ws.Range("K122").Value = 10462.8258
$ws.Range("L122").Value = 1689474
$ws.Range("M122").Value = -8012.825800000001
$ws.Range("N122").Value = -1694374
$ws.Range("H126").Value = 3390.3635
$ws.Range("I126").Value = 7000
$ws.Range("K126").Value = 21000
$ws.Range("M126").Value = -18530
$ws.Range("H132").Value = 6175685.5
$ws.Range("I132").Value = 8335411
$ws.Range("K132").Value = 25006233
$ws.Range("M132").Value = -25003703
$ws.Range("H136").Value = 7610045
$ws.Range("I136").Value = 1611417.1
$ws.Range("J136").Value = 40002636
$ws.Range("K136").Value = 4834251.300000001
$ws.Range("L136").Value = 120007908
$ws.Range("M136").Value = -4831701.300000001
$ws.Range("N136").Value = -120013008

